$d = $word.ActiveDocument

$replacements = @(
    @("78÷9=8, 6", "25÷8=3, 1"),
    @("41÷5=8, 1", "54÷2=27, 0"),
    @("55÷4=13, 3", "34÷2=17, 0"),
    @("85÷9=9, 4", "99÷7=14, 1"),
    @("22÷7=3, 1", "41÷4=10, 1"),
    @("59÷9=6, 5", "33÷7=4, 5"),
    @("84÷9=9, 3", "93÷3=31, 0"),
    @("91÷3=30, 1", "21÷7=3, 0"),
    @("33÷3=11, 0", "59÷3=19, 2"),
    @("58÷7=8, 2", "79÷3=26, 1"),
    @("95÷9=10, 5", "73÷4=18, 1"),
    @("95÷6=15, 5", "78÷5=15, 3"),
    @("25÷5=5, 0", "23÷9=2, 5"),
    @("49÷8=6, 1", "70÷9=7, 7"),
    @("85÷8=10, 5", "27÷3=9, 0"),
    @("65÷4=16, 1", "55÷7=7, 6"),
    @("86÷5=17, 1", "81÷3=27, 0"),
    @("78÷8=9, 6", "53÷8=6, 5"),
    @("50÷9=5, 5", "37÷7=5, 2"),
    @("79÷6=13, 1", "99÷5=19, 4"),
    @("61÷8=7, 5", "11÷3=3, 2"),
    @("23÷6=3, 5", "47÷8=5, 7"),
    @("44÷6=7, 2", "59÷7=8, 3"),
    @("59÷4=14, 3", "89÷6=14, 5"),
    @("11÷2=5, 1", "61÷7=8, 5")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
